$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.007.04"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "1.844.08"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.14"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4670"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3626"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07155"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9160"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07709"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").Value = "1.842.88"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.284"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.421"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.34"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.75%  "
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008594"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "27.033.03"
$ws.Range("E20").Value = "  +2.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.37"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("E22").Value = "  +1.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.64"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.930"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.74"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.28"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.051"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.18"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.67%  "
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08859"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.190"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.874"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.175"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7478"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.470"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.082"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.982"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.02%  "
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05167"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5176"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.912"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1511"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.156"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.48"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4701"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.74"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.604"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.92%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06042"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.59"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.17"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.91%  "
